$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old "py_view" column (C), shifting
# existing py_view..matlab_download columns from C:H to F:K.
$ws.Range("C1:E1").EntireColumn.Insert()

# Give the three freshly inserted columns (C, D, E) a plain, non-bestFit
# width of 11 characters (matches width="11" in the target column defs).
$ws.Columns.Item(3).ColumnWidth = 10.14
$ws.Columns.Item(4).ColumnWidth = 10.14
$ws.Columns.Item(5).ColumnWidth = 10.14

# The new "previews"/"desc" helper column (L) ends up the same width as
# its neighbor (K, formerly r_download/matlab_download) in the target.
$ws.Columns.Item(12).ColumnWidth = 15.17

# Fill the "desc"/"previews" helper columns (L, M) first, row by row,
# matching the order the workbook's shared-string table was built in.
$ws.Range("L1").Value = "desc"
$ws.Range("M1").Value = "previews"
$ws.Range("L2").Value = "Description here."
$ws.Range("M2").Value = "previews/lon_180_to_360/p1.png"
$ws.Range("L3").Value = "Description here."
$ws.Range("M3").Value = "previews/compare_sensors/p1.png"

# Then fill the new modal_id / description / preview_image columns
# (C, D, E), column by column.
$ws.Range("C1").Value = "modal_id"
$ws.Range("C2").Value = "modal-lon"
$ws.Range("C3").Value = "modal-chla"

$ws.Range("D1").Value = "description"
$ws.Range("D2").Value = "Description here"
$ws.Range("D3").Value = "Description here"

$ws.Range("E1").Value = "preview_image"
$ws.Range("E2").Value = "previews/lon.png"
$ws.Range("E3").Value = "previews/chla.png"

# Restore original selection-like active cell (table authoring artifact).
$ws.Range("E6").Select()
